$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$ws.Range("B5").Value = "CDABinaryDataEncoding"
$ws.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$ws.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"
